$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Barco")
$wb.Names.Add('Page1', '=Barco!$B$1:$M$66')
$wb.Names.Add('Page2', '=Barco!$B$67:$M$132')
Write-Host "done"
